$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = -0.1675053202832929
$ws.Range("C2").Value = 0.4571627628656621
$ws.Range("D2").Value = 0.3201243589103758
$ws.Range("E2").Value = 0.5657953330581437
$ws.Range("F2").Value = 0.5608323877005253

# Row 3 (Q1)
$ws.Range("B3").Value = -0.05940777093550652
$ws.Range("C3").Value = 0.470041938014426
$ws.Range("D3").Value = 0.4000049546014873
$ws.Range("E3").Value = 0.6324594489779461
$ws.Range("F3").Value = 0.6553741760501848

# Row 4 (Q2)
$ws.Range("B4").Value = -0.04558429807018168
$ws.Range("C4").Value = 0.6717506146859975
$ws.Range("D4").Value = 0.691874271443171
$ws.Range("E4").Value = 0.8317898000355445
$ws.Range("F4").Value = 0.8674705192030993

# Row 5 (Q3)
$ws.Range("B5").Value = -0.02022650163097881
$ws.Range("C5").Value = 0.7137684849070793
$ws.Range("D5").Value = 0.7655280142584943
$ws.Range("E5").Value = 0.8749445778210722
$ws.Range("F5").Value = 0.9174043782211272

# Row 6 (Q4)
$ws.Range("B6").Value = -0.07724971183834548
$ws.Range("C6").Value = 0.7230381564818515
$ws.Range("D6").Value = 0.7821964137044809
$ws.Range("E6").Value = 0.8844186868810953
$ws.Range("F6").Value = 0.9286961563428431

# Row 7 (Q5)
$ws.Range("B7").Value = -0.1830193523763199
$ws.Range("C7").Value = 0.6659309687502694
$ws.Range("D7").Value = 0.5232006682981125
$ws.Range("E7").Value = 0.7233261147629834
$ws.Range("F7").Value = 0.742238275807101
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = -0.2095810606696748
$ws.Range("C8").Value = 0.830624845422039
$ws.Range("D8").Value = 0.9680408355488034
$ws.Range("E8").Value = 0.9838906623953717
$ws.Range("F8").Value = 1.053062171701582
$ws.Range("G8").Value = 6

# Row 9 (Q7)
$ws.Range("B9").Value = -0.8159383291657744
$ws.Range("C9").Value = 0.8519594231031956
$ws.Range("D9").Value = 1.46036401875364
$ws.Range("E9").Value = 1.208455220003472
$ws.Range("F9").Value = 1.091747678095863
$ws.Range("G9").Value = 3

# Row 10 (new, Q8) - copy formatting from row 9's A cell so the new
# label cell picks up the same style (bold, bordered, centered) as the
# other row labels.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.9652956092697305
$ws.Range("C10").Value = 0.9652956092697305
$ws.Range("D10").Value = 0.9317956132754202
$ws.Range("E10").Value = 0.9652956092697305
$ws.Range("G10").Value = 1
